# Backlog.xlsx — "Nu lämnar jag in sprinten"
#
# Sheet "Blad1" columns: A=ID, B=Uppgift, C=Typ, D=Prioritering,
# E=Status (conditional "Bra"/"Neutral"/"Dålig" cell styles), F=Sprint,
# G=Kommentar.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 15-16 ("Lägg till shoppingvagn" / "Lägg till köpfunktionalitet") are
# now done: flip their Status cell style from "Neutral" to "Bra" (green)
# and clear out the now-stale Kommentar notes.
$ws.Range("E15").Style = "Bra"
$ws.Range("G15").ClearContents()

$ws.Range("E16").Style = "Bra"
$ws.Range("G16").ClearContents()

# Row 22 ("Skydd mot SQL injection") is re-typed as a System task, and its
# priority (along with row 23 "Skydd mot session hijacking") is bumped down
# to "Väldigt låg".
$ws.Range("C22").Value = "System"
$ws.Range("D22").Value = "Väldigt låg"
$ws.Range("D23").Value = "Väldigt låg"

# Leave the cursor where the author last left it.
$ws.Range("D24").Select() | Out-Null
